# Applies the "han added his part" edit: appends two new analysis
# paragraphs (separated by blank paragraphs) after the existing
# "...comparatively similar." paragraph, moves the _GoBack bookmark to
# the end of the new last paragraph, and leaves two trailing blank
# paragraphs at the end of the document body.

$d = $word.ActiveDocument

# The paragraph we are extending is the last paragraph in the document
# (the one ending "...comparatively similar." which also currently
# carries the _GoBack bookmark).
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $lastPara.Range

$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="1658FCCC" w14:textId="2F8AB7FC" w:rsidR="00D351E5" w:rsidRDefault="00D351E5" w:rsidP="00D351E5"><w:r><w:t xml:space="preserve">From our data analysis, we found that highest number of Data Science jobs were available in California, followed by New York, Virginia, Texas and then Massachusetts. Relatively high number of Data Science jobs were available in California, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>where as</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> in Texas, the availability all these job types were comparatively similar.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">For a specific skill, the number of jobs is significantly different among 3 job types. </w:t></w:r><w:r><w:t>Python and SQL are essential skills required by all 3 job types. R is needed for both data analyst and data scientist which also requires tableau and machine learning respectively. Whereas, Data engineer needs Hadoop, spark and java.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">For a specific skill, the number of jobs is significantly different </w:t></w:r><w:r><w:t>in salary categories</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r><w:t xml:space="preserve"> The need for python, spark, java, Hadoop and machine learning increases as the salary increasing; while the need for SQL, r, tableau and SAS stays stable across all salary categories.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p/><w:p/>
'@

# Replacing the whole paragraph range with OOXML that reproduces the
# original runs (without the bookmark) followed by the new paragraphs;
# the _GoBack bookmark is recreated at the end of the final new
# paragraph, matching where Word leaves it after the last edit.
[void]$r.InsertXML($xml)
